$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet (report title -> "rejected theses" variant)
$ws.Name = "DANH SÁCH ĐỀ TÀI BỊ TỪ CHỐI"

# 2. Insert a new column before the current column E ("Năm thực hiện")
#    to hold the lecturer-name field. Insert shifts F.. onward to G.. etc.
$ws.Columns("E:E").Insert()

# 3. Clone the formatting that column F (the column right after the newly
#    inserted, blank column E) carries for the header band / value band
#    rows so the new column matches the rest of the table visually.
$ws.Range("F1:F3").Copy()
$ws.Range("E1:E3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4. Put in the new header / placeholder text
$ws.Range("E2").Value = "Giảng viên ra đề"
$ws.Range("E3").Value = "{{Items.LecturerName}}"

# 5. Re-apply the column widths to match the target layout. The insert
#    operation shifts the pre-existing per-column widths one slot to the
#    right along with the cells, so every column from F onward needs its
#    width restored to what it was (by letter) before the insert, and the
#    two new columns (E, K) get their own widths.
$ws.Columns("E").ColumnWidth = 29.333333333333332
$ws.Columns("F").ColumnWidth = 20.166666666666668
$ws.Columns("G").ColumnWidth = 19.333333333333332
$ws.Columns("H").ColumnWidth = 23
$ws.Columns("I").ColumnWidth = 28.333333333333332
$ws.Columns("J").ColumnWidth = 47
$ws.Columns("K").ColumnWidth = 34.333333333333336

# 6. Update the active selection to the new lecturer-name cell
$ws.Range("E3").Select()
